# Append 12 more rows (A=204..215) of data to the "월_중국연휴수" normalized
# sheet, continuing the existing A/B series that ran through row 205
# (A205=203). New data goes into rows 206-217.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @(204, 0.3230769230769231),
    @(205, 0.3118881118881119),
    @(206, 0.4461538461538461),
    @(207, 0.4461538461538461),
    @(208, 0.4461538461538461),
    @(209, 0.4461538461538461),
    @(210, 0.9384615384615382),
    @(211, 0.4461538461538461),
    @(212, 0.4461538461538461),
    @(213, 0.9384615384615382),
    @(214, 0.4461538461538461),
    @(215, 0.4461538461538461)
)

$startRow = 206

# Template cell for the formatting already used on column A (bold, boxed,
# centered/top-aligned) so the appended rows match the existing style.
$templateA = $ws.Range("A205")

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $aCell = $ws.Cells.Item($row, 1)
    $bCell = $ws.Cells.Item($row, 2)

    $aCell.Value = $newData[$i][0]
    $bCell.Value = $newData[$i][1]

    $templateA.Copy() | Out-Null
    $aCell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

$excel.CutCopyMode = 0
